# Apply the cryptos-list price/volume refresh for Sun Apr  2 17:20:49 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.201.76'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.801.55'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.35'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5282'
$ws.Range("E7").Value = '  +3.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3825'
$ws.Range("E8").Value = '  -2.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08012'
$ws.Range("E9").Value = '  +3.07%  '
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.100'
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.327'
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").Value = '1.806.17'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.330'
$ws.Range("E16").Value = '  -1.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.97'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001098'
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06601'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.36'
$ws.Range("E21").Value = '  -2.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.969'
$ws.Range("E22").Value = '  -1.85%  '
$ws.Range("D23").Value = '28.239.44'
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.21'
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.232'
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.38'
$ws.Range("E26").Value = '  +3.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.53'
$ws.Range("E27").Value = '  -2.63%  '
$ws.Range("D28").Value = '2.007.81'
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.26'
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1092'
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.061'
$ws.Range("E32").Value = '  -3.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.657'
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.554'
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07311'
$ws.Range("E35").Value = '  +3.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.26'
$ws.Range("E36").Value = '  +9.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.905'
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2168'
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02312'
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.076'
$ws.Range("E40").Value = '  -2.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6215'
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.165'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.27'
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6006'
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.760'
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '126.82'
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.208'
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.928'
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06830'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.16'
$ws.Range("E51").Value = '  -1.12%  '
